# feat: add 2022-Q3 data
#
# The workbook currently has two sheets: "总计" (totals) and "2022-Q2"
# (fund holdings for 2022-Q2). This script:
#   1. Duplicates the "2022-Q2" sheet (the duplicate keeps the old data
#      and stays named "2022-Q2"), placing the duplicate right after the
#      original.
#   2. Renames/repurposes the original "2022-Q2" sheet into "2022-Q3" and
#      replaces its contents with the new Q3 fund-holdings data.
#   3. Inserts a new summary row into "总计" for 2022-Q3 (above the
#      existing 2022-Q2 summary row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate "2022-Q2" so we keep its data around under its own
# tab, then repurpose the original tab as "2022-Q3".
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($null, $wsQ2)
$wsQ2Copy = $wb.Worksheets.Item($wsQ2.Index + 1)

$wsQ2.Name = "2022-Q3"
$wsQ2Copy.Name = "2022-Q2"

$wsQ3 = $wsQ2

# ---------------------------------------------------------------------
# Step 2: overwrite the (former 2022-Q2, now 2022-Q3) sheet with the new
# quarter's fund holdings. Header row (row 1) stays the same; only the
# data rows change content, and the row count grows from 1 to 5 data
# rows.
# ---------------------------------------------------------------------

# Keep numbers like "3.80"/"0.01"/"001668" as text so the trailing and
# leading zeros survive (matches how the source data is stored).
$wsQ3.Range("B2:B6").NumberFormat = "@"
$wsQ3.Range("D2:G6").NumberFormat = "@"

$wsQ3.Range("A2").Value() = 0
$wsQ3.Range("B2").Value() = "513360"
$wsQ3.Range("C2").Value() = "博时中证全球中国教育主题ETF（QDII）"
$wsQ3.Range("D2").Value() = "4.81"
$wsQ3.Range("E2").Value() = "99.43"
$wsQ3.Range("F2").Value() = "24.38"
$wsQ3.Range("G2").Value() = "1.1727"
$wsQ3.Range("H2").Value() = 1

$wsQ3.Range("A3").Value() = 1
$wsQ3.Range("B3").Value() = "001668"
$wsQ3.Range("C3").Value() = "汇添富全球移动互联灵活配置混合（QDII）A"
$wsQ3.Range("D3").Value() = "12.06"
$wsQ3.Range("E3").Value() = "90.88"
$wsQ3.Range("F3").Value() = "3.11"
$wsQ3.Range("G3").Value() = "0.3751"
$wsQ3.Range("H3").Value() = 5

$wsQ3.Range("A4").Value() = 2
$wsQ3.Range("B4").Value() = "457001"
$wsQ3.Range("C4").Value() = "国富亚洲机会股票（QDII）"
$wsQ3.Range("D4").Value() = "3.80"
$wsQ3.Range("E4").Value() = "83.80"
$wsQ3.Range("F4").Value() = "3.03"
$wsQ3.Range("G4").Value() = "0.1151"
$wsQ3.Range("H4").Value() = 8

$wsQ3.Range("A5").Value() = 3
$wsQ3.Range("B5").Value() = "015203"
$wsQ3.Range("C5").Value() = "汇添富全球移动互联灵活配置混合（QDII）D"
$wsQ3.Range("D5").Value() = "0.04"
$wsQ3.Range("E5").Value() = "90.88"
$wsQ3.Range("F5").Value() = "3.11"
$wsQ3.Range("G5").Value() = "0.0012"
$wsQ3.Range("H5").Value() = 5

$wsQ3.Range("A6").Value() = 4
$wsQ3.Range("B6").Value() = "015202"
$wsQ3.Range("C6").Value() = "汇添富全球移动互联灵活配置混合（QDII）C"
$wsQ3.Range("D6").Value() = "0.01"
$wsQ3.Range("E6").Value() = "90.88"
$wsQ3.Range("F6").Value() = "3.11"
$wsQ3.Range("G6").Value() = "0.0003"
$wsQ3.Range("H6").Value() = 5

# ---------------------------------------------------------------------
# Step 3: insert a 2022-Q3 row at the top of the "总计" summary sheet,
# pushing the existing 2022-Q2 row down.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))

$wsTotal.Range("A2").Value() = 0
$wsTotal.Range("B2").Value() = "2022-Q3"
$wsTotal.Range("C2").Value() = 5
$wsTotal.Range("D2").Value() = 1.66

$wsTotal.Range("A3").Value() = 1
